$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1780821917808219
$ws.Range("C2").Value = 0.571917808219178
$ws.Range("J2").Value = 0.03082191780821918
$ws.Range("O2").Value = 0.003424657534246575
$ws.Range("P2").Value = 0.1506849315068493
$ws.Range("S2").Value = 0.06506849315068493
$ws.Range("B3").Value = 0.005555555555555556
$ws.Range("C3").Value = 0.01111111111111111
$ws.Range("J3").Value = 0.02777777777777778
$ws.Range("P3").Value = 0.7722222222222223
$ws.Range("S3").Value = 0.1833333333333333
$ws.Range("J4").Value = 0.0625
$ws.Range("P4").Value = 0.6041666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("B6").Value = 0.05527638190954774
$ws.Range("D6").Value = 0.01507537688442211
$ws.Range("F6").Value = 0.04522613065326633
$ws.Range("J6").Value = 0.185929648241206
$ws.Range("O6").Value = 0.02512562814070352
$ws.Range("Q6").Value = 0.1557788944723618
$ws.Range("R6").Value = 0.1005025125628141
$ws.Range("S6").Value = 0.4170854271356784
$ws.Range("B7").Value = 0.1106870229007634
$ws.Range("D7").Value = 0.007633587786259542
$ws.Range("F7").Value = 0.05343511450381679
$ws.Range("J7").Value = 0.1412213740458015
$ws.Range("O7").Value = 0.03053435114503817
$ws.Range("Q7").Value = 0.1755725190839695
$ws.Range("R7").Value = 0.1068702290076336
$ws.Range("S7").Value = 0.3740458015267176
$ws.Range("B8").Value = 0.07520325203252033
$ws.Range("D8").Value = 0.02439024390243903
$ws.Range("F8").Value = 0.03861788617886179
$ws.Range("J8").Value = 0.1239837398373984
$ws.Range("O8").Value = 0.01016260162601626
$ws.Range("Q8").Value = 0.1483739837398374
$ws.Range("R8").Value = 0.1158536585365854
$ws.Range("S8").Value = 0.4634146341463415
$ws.Range("B9").Value = 0.05504587155963303
$ws.Range("D9").Value = 0.02293577981651376
$ws.Range("E9").Value = 0.004587155963302753
$ws.Range("F9").Value = 0.05045871559633028
$ws.Range("J9").Value = 0.1146788990825688
$ws.Range("O9").Value = 0.009174311926605505
$ws.Range("Q9").Value = 0.1330275229357798
$ws.Range("R9").Value = 0.1284403669724771
$ws.Range("S9").Value = 0.481651376146789
$ws.Range("B10").Value = 0.1038575667655786
$ws.Range("D10").Value = 0.01928783382789317
$ws.Range("F10").Value = 0.06973293768545995
$ws.Range("J10").Value = 0.120919881305638
$ws.Range("O10").Value = 0.02596439169139466
$ws.Range("Q10").Value = 0.1862017804154303
$ws.Range("R10").Value = 0.1001483679525223
$ws.Range("S10").Value = 0.3738872403560831
$ws.Range("G11").Value = 0.1503416856492027
$ws.Range("J11").Value = 0.09111617312072894
$ws.Range("K11").Value = 0.2072892938496583
$ws.Range("L11").Value = 0.5193621867881549
$ws.Range("S11").Value = 0.03189066059225513
$ws.Range("G12").Value = 0.7100840336134454
$ws.Range("J12").Value = 0.2100840336134454
$ws.Range("K12").Value = 0.004201680672268907
$ws.Range("L12").Value = 0.03781512605042017
$ws.Range("S12").Value = 0.03781512605042017
$ws.Range("F13").Value = 0.02040816326530612
$ws.Range("G13").Value = 0.6530612244897959
$ws.Range("J13").Value = 0.2653061224489796
$ws.Range("S13").Value = 0.06122448979591837
$ws.Range("F15").Value = 0.01659751037344398
$ws.Range("H15").Value = 0.1784232365145228
$ws.Range("I15").Value = 0.05394190871369295
$ws.Range("J15").Value = 0.3319502074688797
$ws.Range("K15").Value = 0.07053941908713693
$ws.Range("M15").Value = 0.01659751037344398
$ws.Range("O15").Value = 0.06224066390041494
$ws.Range("S15").Value = 0.2697095435684647
$ws.Range("F16").Value = 0.004926108374384237
$ws.Range("H16").Value = 0.1428571428571428
$ws.Range("I16").Value = 0.09852216748768473
$ws.Range("J16").Value = 0.3891625615763547
$ws.Range("K16").Value = 0.1477832512315271
$ws.Range("M16").Value = 0.004926108374384237
$ws.Range("N16").Value = 0.009852216748768473
$ws.Range("O16").Value = 0.02463054187192118
$ws.Range("S16").Value = 0.1773399014778325
$ws.Range("F17").Value = 0.01157407407407407
$ws.Range("H17").Value = 0.2060185185185185
$ws.Range("I17").Value = 0.09490740740740741
$ws.Range("J17").Value = 0.349537037037037
$ws.Range("K17").Value = 0.1157407407407407
$ws.Range("M17").Value = 0.02083333333333333
$ws.Range("O17").Value = 0.06712962962962964
$ws.Range("S17").Value = 0.1342592592592593
$ws.Range("F18").Value = 0.01465201465201465
$ws.Range("H18").Value = 0.1538461538461539
$ws.Range("I18").Value = 0.1172161172161172
$ws.Range("J18").Value = 0.3882783882783883
$ws.Range("K18").Value = 0.1501831501831502
$ws.Range("M18").Value = 0.01831501831501832
$ws.Range("O18").Value = 0.0695970695970696
$ws.Range("S18").Value = 0.08791208791208792
$ws.Range("F19").Value = 0.01147873058744092
$ws.Range("H19").Value = 0.1971640783254558
$ws.Range("I19").Value = 0.07900067521944631
$ws.Range("J19").Value = 0.3517893315327482
$ws.Range("K19").Value = 0.1411208642808913
$ws.Range("M19").Value = 0.02093180283592167
$ws.Range("N19").Value = 0.0006752194463200541
$ws.Range("O19").Value = 0.05671843349088454
$ws.Range("S19").Value = 0.1411208642808913
